# "Ran Prod Verification Script"
# The verification script stamps a fresh "DateProd" timestamp (column B) into
# row 2 (and rows 3-5, for the 4-row "-Generic" sheets) of every "-Prod"
# result sheet each time it runs. This reproduces the timestamps written by
# that run.

$wb = $excel.ActiveWorkbook

function Set-ProdDate {
    param([string]$SheetName, [hashtable]$CellValues)

    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($cellRef in $CellValues.Keys) {
        $ws.Range($cellRef).Value = $CellValues[$cellRef]
    }
}

Set-ProdDate "VT-C-DebitCredit-DualCF-Generic" @{
    "B2" = "Thu Aug 28 08:15:12 IST 2025"
}

Set-ProdDate "VT-C-DebitCredit-NoCF-Generic" @{
    "B2" = "Thu Aug 28 08:16:15 IST 2025"
}

Set-ProdDate "VT-C-DebitCredit-SingleCF-Gener" @{
    "B2" = "Thu Aug 28 08:17:17 IST 2025"
}

Set-ProdDate "VT-P-DebitCredit-DualCF-Generic" @{
    "B2" = "Thu Aug 28 08:18:19 IST 2025"
    "B3" = "Thu Aug 28 08:19:28 IST 2025"
    "B4" = "Thu Aug 28 08:20:35 IST 2025"
    "B5" = "Thu Aug 28 08:21:42 IST 2025"
}

Set-ProdDate "VT-P-DebitCredit-NoCF-Generic" @{
    "B2" = "Thu Aug 28 08:22:45 IST 2025"
    "B3" = "Thu Aug 28 08:23:48 IST 2025"
    "B4" = "Thu Aug 28 08:24:53 IST 2025"
    "B5" = "Thu Aug 28 08:26:01 IST 2025"
}

Set-ProdDate "VT-P-DebitCredit-SingleCF-Gener" @{
    "B2" = "Thu Aug 28 08:27:07 IST 2025"
    "B3" = "Thu Aug 28 08:28:11 IST 2025"
    "B4" = "Thu Aug 28 08:29:16 IST 2025"
    "B5" = "Thu Aug 28 08:30:26 IST 2025"
}

Set-ProdDate "VT-P-DebitVoid-DualCF-Generic" @{
    "B2" = "Thu Aug 28 08:31:33 IST 2025"
}

Set-ProdDate "VT-P-DebitVoid-NoCF-Generic" @{
    "B2" = "Thu Aug 28 08:32:47 IST 2025"
    "B3" = "Thu Aug 28 08:33:54 IST 2025"
    "B4" = "Thu Aug 28 08:34:59 IST 2025"
    "B5" = "Thu Aug 28 08:36:01 IST 2025"
}

Set-ProdDate "VT-P-DebitVoid-SingleCF-Generic" @{
    "B2" = "Thu Aug 28 08:37:04 IST 2025"
    "B3" = "Thu Aug 28 08:38:13 IST 2025"
    "B4" = "Thu Aug 28 08:39:23 IST 2025"
    "B5" = "Thu Aug 28 08:40:30 IST 2025"
}
